$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for all data rows (2..108)
# from serial date 45207 (2023-10-08) to 45208 (2023-10-09).
for ($row = 2; $row -le 108; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45207) {
        $cell.Value2 = 45208
    }
}

# Rename folder references in hyperlink formulas for rows 2 and 3:
#   Logging_VASTERVIK -> Logging_0883
#   Logging_MONSTERAS -> Logging_0861
$row2Cols = @("S", "T", "U", "V", "W", "X", "Y")
foreach ($col in $row2Cols) {
    $cell = $ws.Range($col + "2")
    $cell.Formula = $cell.Formula -replace "Logging_VASTERVIK", "Logging_0883"
}

foreach ($col in $row2Cols) {
    $cell = $ws.Range($col + "3")
    $cell.Formula = $cell.Formula -replace "Logging_MONSTERAS", "Logging_0861"
}
